$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so values like "325.73" are not
# auto-converted to numbers by Excel (matches source inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "49.897.98"
$ws.Range("E2").Value = "  +4.51%  "

$ws.Range("D3").Value = "2.674.33"
$ws.Range("E3").Value = "  +7.91%  "

$ws.Range("D4").Value = "0.999"

$ws.Range("D5").Value = "113.59"
$ws.Range("E5").Value = "  +9.48%  "

$ws.Range("D6").Value = "325.73"
$ws.Range("E6").Value = "  +3.08%  "

$ws.Range("D7").Value = "0.528"
$ws.Range("E7").Value = "  +2.25%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "0.552"
$ws.Range("E9").Value = "  +3.76%  "

$ws.Range("D10").Value = "40.91"
$ws.Range("E10").Value = "  +6.23%  "

$ws.Range("D11").Value = "20.12"
$ws.Range("E11").Value = "  -1.54%  "

$ws.Range("D12").Value = "0.0823"
$ws.Range("E12").Value = "  +3.43%  "

$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "7.36"
$ws.Range("E14").Value = "  +5.03%  "

$ws.Range("D15").Value = "3.082.52"
$ws.Range("E15").Value = "  +7.58%  "

$ws.Range("D16").Value = "2.665.15"
$ws.Range("E16").Value = "  +7.31%  "

$ws.Range("D17").Value = "0.874"
$ws.Range("E17").Value = "  +6.55%  "

$ws.Range("D18").Value = "49.823.58"
$ws.Range("E18").Value = "  +4.52%  "

$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  +4.19%  "

$ws.Range("D20").Value = "6.79"
$ws.Range("E20").Value = "  +4.43%  "

$ws.Range("D21").Value = "2.90"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("D22").Value = "0.0₃0958"

$ws.Range("D23").Value = "278.28"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "71.78"
$ws.Range("E24").Value = "  +1.62%  "

$ws.Range("E25").Value = "  +3.46%  "

$ws.Range("D26").Value = "26.84"
$ws.Range("E26").Value = "  +4.83%  "

$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +6.21%  "

$ws.Range("E29").Value = "  +1.27%  "

$ws.Range("D30").Value = "36.23"
$ws.Range("E30").Value = "  +6.02%  "

$ws.Range("D31").Value = "0.141"
$ws.Range("E31").Value = "  +3.96%  "

$ws.Range("D32").Value = "50.22"
$ws.Range("E32").Value = "  +2.21%  "

$ws.Range("E33").Value = "  +4.74%  "

$ws.Range("D34").Value = "19.51"
$ws.Range("E34").Value = "  +3.55%  "

$ws.Range("E35").Value = "  +5.93%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +13.72%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D39").Value = "3.17"
$ws.Range("E39").Value = "  +11.70%  "

$ws.Range("D40").Value = "125.02"
$ws.Range("E40").Value = "  +1.81%  "

$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").Value = "22.75"
$ws.Range("E42").Value = "  +5.12%  "

$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("E44").Value = "  +6.73%  "

$ws.Range("D45").Value = "2.107.49"
$ws.Range("E45").Value = "  +6.03%  "

$ws.Range("E46").Value = "  +6.33%  "

$ws.Range("E47").Value = "  +15.59%  "

$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +8.93%  "

$ws.Range("D49").Value = "9.04"
$ws.Range("E49").Value = "  +1.93%  "

$ws.Range("D50").Value = "5.36"
$ws.Range("E50").Value = "  +6.35%  "

$ws.Range("D51").Value = "59.55"
$ws.Range("E51").Value = "  +7.50%  "
